# ---------------------------------------------------------------------------
# Edit script: "Diving Deep into Cyber Resilience" -> "The Enchanting
# Symphony: Unveiling the Beauty of Chemistry" (Calibri font set, doc 0135)
# ---------------------------------------------------------------------------
$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Host "WARNING: not found -> [$old]"
    }
}

# --- Title -------------------------------------------------------------
Replace-Text "Diving Deep into Cyber Resilience: A Cybersecurity Imperative" "The Enchanting Symphony: Unveiling the Beauty of Chemistry"

# --- Byline (author name) ----------------------------------------------
Replace-Text "Sophia Henderson" "Dr. Elise Campbell"

# --- Email paragraph -----------------------------------------------------
Replace-Text "sophia" "ecampbell@highschoolofchemistry"
Replace-Text "henderson@emailuniversity.edu" "org"

# --- Body paragraph 1 (four sentences before first <w:br/><w:br/>) ------
Replace-Text "In a world where the reliance on digital technologies continues to grow exponentially, the significance of cybersecurity has never been more apparent" "Every element of chemistry reflects the harmony of life"
Replace-Text "As the realm of cyberspace expands, so do the threats that lurk within it, challenging organizations and individuals alike" "In this extraordinary symphony of matter, the blend of elements creates a magnificent tapestry that captivates the senses"
Replace-Text "Embracing cyber resilience is a critical imperative, a fundamental shift in mindset that necessitates proactive steps to safeguard valuable assets" "The interconnectedness between the periodic table, the elements, and the intricate structures of molecules orchestrates a fascinating composition"
Replace-Text "This comprehensive strategy involves implementing robust safeguards, fostering a security-conscious culture, and ensuring swift recovery from potential incidents" "As you embark on this scientific journey, you will become a chemist, an artist, and a composer, blending ingredients, conducting experiments, and unraveling the enigma of the chemical world"

# --- Second block (after first <w:br/><w:br/>) ---------------------------
Replace-Text "As malicious actors continue to devise sophisticated attacks, traditional defensive strategies often prove inadequate" "At the heart of chemistry is the enigma of atoms"
Replace-Text "Cyber resilience requires organizations to adopt an agile and proactive approach, continuously adapting to evolving threats and remaining prepared for contingencies" "These subatomic particles, like tiny cosmic dancers, weave an intricate ballet of electrons, protons, and neutrons"
Replace-Text "This proactive stance entails implementing rigorous security measures such as multi-factor authentication, encryption, and regular patching of systems" "As they waltz around the nucleus, energy radiates, orchestrating the symphony of particles that form molecules, the building blocks of all matter"
Replace-Text "Additionally, promoting a culture of cybersecurity awareness among employees and stakeholders empowers them to recognize and mitigate potential threats" "We will explore the periodic table, a symphony of elements arranged in a profound and mystical order, revealing the symphony of qualities and characteristics that each possesses. Every atom, an individual instrument, plays its role in shaping the melody of matter"

# --- Third block (after second <w:br/><w:br/>) ----------------------------
Replace-Text "Furthermore, the ability to recover swiftly from a cyber incident is a crucial aspect of cyber resilience" "Finally, we delve into the captivating realm of chemical reactions, the enigmatic dance that takes place when atoms rearrange themselves"
Replace-Text "With prompt detection mechanisms in place, organizations can respond rapidly to contain and eradicate threats" "Sparks fly, colors change, and new substances are formed"
Replace-Text "Regularly testing incident response plans, conducting security audits, and partnering with cybersecurity experts are essential steps in ensuring swift recovery and minimizing potential damages" "From the fizz of a baking soda volcano to the glow of a firecracker, these reactions ignite a symphony of senses"
Replace-Text "These recovery efforts go beyond technological measures, encompassing thorough communication with stakeholders and clients, mitigating reputational damage, and safeguarding customer trust" "We will grapple with the energies that bind and break molecules, understanding the intricate balance that guides these chemical transformations, and hear the chorus of molecules singing in perfect pitch"

# --- Summary paragraph -----------------------------------------------------
Replace-Text "In today's digitally connected world, cyber resilience has emerged as a critical imperative for organizations and individuals alike" "Here, you will venture through the captivating universe of chemistry, exploring the harmony of elements, the compounds they form, and the reactions that transform our world"
Replace-Text "Embracing a proactive approach, implementing robust security measures, cultivating a security-conscious culture, and ensuring rapid recovery from incidents are fundamental aspects of effective cyber resilience" "With each step, you will discover the beauty of science hidden within the symphony of matter"
Replace-Text "By doing so, organizations can significantly reduce the impact of cyber threats, protect valuable assets, maintain seamless operations, and uphold customer trust in the face of evolving cybersecurity challenges" "Chemistry is not just a collection of facts and figures; it is an art form, a dance of particles, and a breathtaking display of nature's boundless creativity. Prepare to be captivated by this enchanting symphony and fall in love with the allure of chemistry"

# --- Trailing empty paragraph ----------------------------------------------
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()
